$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21 (shifts existing rows 21-140 down to 22-141)
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 with the new data record.
$ws.Cells.Item(21, 1).Value = 4
$ws.Cells.Item(21, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(21, 3).Value = "Los Lagos"
$ws.Cells.Item(21, 4).Value = 44819
$ws.Cells.Item(21, 5).Value = 10
$ws.Cells.Item(21, 6).Value = 100112052
$ws.Cells.Item(21, 7).Value = "Albahaca"
$ws.Cells.Item(21, 8).Value = "Sin especificar"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 80
$ws.Cells.Item(21, 11).Value = 6500
$ws.Cells.Item(21, 12).Value = 6500
$ws.Cells.Item(21, 13).Value = 6500
$ws.Cells.Item(21, 14).Value = "$/paquete"
$ws.Cells.Item(21, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(21, 16).Value = 6500
$ws.Cells.Item(21, 17).Value = 1
$ws.Cells.Item(21, 18).Value = "Hortaliza"
